$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 8.5
$ws.Range("AD3").Value = 23
$ws.Range("J3").Value = 1.54
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 2.62
$ws.Range("K4").Value = 1.87
$ws.Range("L4").Value = 6
$ws.Range("M4").Value = 1.1
$ws.Range("O4").Value = 1.54
$ws.Range("R4").Value = 1.41
$ws.Range("T4").Value = 1.1
$ws.Range("W4").Value = 2.5
$ws.Range("X4").Value = 1.5
$ws.Range("Y4").Value = 4.75
$ws.Range("G5").Value = 2.8
$ws.Range("H5").Value = 2.75
$ws.Range("K5").Value = 1.8
$ws.Range("M5").Value = 1.1
$ws.Range("O5").Value = 1.54
$ws.Range("T5").Value = 1.1
$ws.Range("H6").Value = 2.75
$ws.Range("I6").Value = 2.55
$ws.Range("K6").Value = 1.77
$ws.Range("M6").Value = 1.13
$ws.Range("O6").Value = 1.69
$ws.Range("T6").Value = 1.05
$ws.Range("AP7").Value = 6.2
$ws.Range("AR7").Value = 2.49
$ws.Range("AS7").Value = 1.54
$ws.Range("G7").Value = 2.15
$ws.Range("K7").Value = 1.77
$ws.Range("M7").Value = 1.11
$ws.Range("O7").Value = 1.69
$ws.Range("T7").Value = 1.05
$ws.Range("AG8").Value = 21
$ws.Range("AN8").Value = 41
$ws.Range("G8").Value = 1.75
$ws.Range("I8").Value = 5
$ws.Range("M8").Value = 1.05
$ws.Range("O8").Value = 1.41
$ws.Range("P8").Value = 2.62
$ws.Range("T8").Value = 1.13
$ws.Range("Y8").Value = 5.5
$ws.Range("G9").Value = 2.2
$ws.Range("H9").Value = 2.88
$ws.Range("K9").Value = 1.8
$ws.Range("M9").Value = 1.1
$ws.Range("O9").Value = 1.54
$ws.Range("T9").Value = 1.1
$ws.Range("G10").Value = 1.9
$ws.Range("H10").Value = 2.9
$ws.Range("R10").Value = 1.3
$ws.Range("AK14").Value = 19
$ws.Range("AM14").Value = 41
$ws.Range("AN14").Value = 29
$ws.Range("G14").Value = 2.05
$ws.Range("I14").Value = 3.5
$ws.Range("AA15").Value = 9
$ws.Range("AB15").Value = 17
$ws.Range("AK15").Value = 19
$ws.Range("AN15").Value = 29
$ws.Range("G15").Value = 1.95
$ws.Range("I15").Value = 3.8
$ws.Range("L15").Value = 4.33
$ws.Range("Y15").Value = 7.5
$ws.Range("AF17").Value = 5.5
$ws.Range("K17").Value = 1.87
$ws.Range("M17").Value = 1.13
$ws.Range("N17").Value = 6
$ws.Range("K18").Value = 1.87
$ws.Range("U18").Value = 1.62
$ws.Range("L19").Value = 2.87
$ws.Range("U19").Value = 1.36
$ws.Range("U20").Value = 1.3
$ws.Range("U24").Value = 1.4
$ws.Range("O29").Value = 1.1
$ws.Range("S29").Value = 1.83
$ws.Range("T29").Value = 1.83
$ws.Range("W29").Value = 1.41
$ws.Range("X29").Value = 2.62
$ws.Range("AK30").Value = 11
$ws.Range("M30").Value = 1.04
$ws.Range("O30").Value = 1.22
$ws.Range("T30").Value = 1.4
$ws.Range("W30").Value = 1.58
$ws.Range("Y30").Value = 12
$ws.Range("AP31").Value = 3.2
$ws.Range("AQ31").Value = 1.35
$ws.Range("AR31").Value = 1.63
$ws.Range("AS31").Value = 2.28
$ws.Range("M31").Value = 1.08
$ws.Range("O31").Value = 1.33
$ws.Range("T31").Value = 1.25
$ws.Range("W31").Value = 1.8
$ws.Range("X31").Value = 1.8
$ws.Range("G32").Value = 1.95
$ws.Range("I32").Value = 3.7
$ws.Range("J32").Value = 2.6
$ws.Range("L32").Value = 4
$ws.Range("M32").Value = 1.04
$ws.Range("N32").Value = 12
$ws.Range("O32").Value = 1.22
$ws.Range("T32").Value = 1.4
$ws.Range("W32").Value = 1.63
$ws.Range("X32").Value = 2.1
$ws.Range("Y32").Value = 8.5
$ws.Range("Z32").Value = 10
$ws.Range("AD33").Value = 41
$ws.Range("AJ33").Value = 7.5
$ws.Range("AM33").Value = 13
$ws.Range("G33").Value = 4.2
$ws.Range("M33").Value = 1.04
$ws.Range("N33").Value = 13
$ws.Range("O33").Value = 1.25
$ws.Range("Q33").Value = 1.85
$ws.Range("R33").Value = 2
$ws.Range("T33").Value = 1.36
$ws.Range("W33").Value = 1.77
$ws.Range("X33").Value = 1.87
$ws.Range("Y33").Value = 12
$ws.Range("K34").Value = 1.95
$ws.Range("M34").Value = 1.1
$ws.Range("O34").Value = 1.5
$ws.Range("R34").Value = 1.47
$ws.Range("T34").Value = 1.17
$ws.Range("X34").Value = 1.67
$ws.Range("M35").Value = 1.06
$ws.Range("O35").Value = 1.36
$ws.Range("R35").Value = 1.63
$ws.Range("T35").Value = 1.22
$ws.Range("W35").Value = 1.87
$ws.Range("X35").Value = 1.77
$ws.Range("AB37").Value = 70
$ws.Range("AD37").Value = 55
$ws.Range("AI37").Value = 900
$ws.Range("AJ37").Value = 5.9
$ws.Range("AK37").Value = 8
$ws.Range("AM37").Value = 16
$ws.Range("G37").Value = 4.1
$ws.Range("I37").Value = 1.87
$ws.Range("J37").Value = 4.6
$ws.Range("P37").Value = 2.72
$ws.Range("R37").Value = 1.6
$ws.Range("S37").Value = 3.8
$ws.Range("U37").Value = 1.45
$ws.Range("V37").Value = 2.55
$ws.Range("W37").Value = 1.98
$ws.Range("Y37").Value = 10.25
$ws.Range("Z37").Value = 23
$ws.Range("W39").Value = 1.77
$ws.Range("X39").Value = 1.87
$ws.Range("AA40").Value = 9.5
$ws.Range("AB40").Value = 30
$ws.Range("AC40").Value = 20
$ws.Range("AE40").Value = 7.3
$ws.Range("AG40").Value = 12.5
$ws.Range("AM40").Value = 29
$ws.Range("AO40").Value = 28
$ws.Range("G40").Value = 2.62
$ws.Range("I40").Value = 2.57
$ws.Range("J40").Value = 3.15
$ws.Range("K40").Value = 2.07
$ws.Range("N40").Value = 7.3
$ws.Range("P40").Value = 3.3
$ws.Range("Q40").Value = 1.85
$ws.Range("R40").Value = 1.85
$ws.Range("S40").Value = 3
$ws.Range("T40").Value = 1.34
$ws.Range("U40").Value = 1.4
$ws.Range("V40").Value = 2.7
$ws.Range("X40").Value = 2.1
$ws.Range("Y40").Value = 9.5
$ws.Range("Z40").Value = 14.5
$ws.Range("M41").Value = 1.01
$ws.Range("O41").Value = 1.11
$ws.Range("T41").Value = 1.63
$ws.Range("M42").Value = 1.03
$ws.Range("O42").Value = 1.19
$ws.Range("T42").Value = 1.37
$ws.Range("AA43").Value = 10
$ws.Range("AB43").Value = 26
$ws.Range("AC43").Value = 21
$ws.Range("AI43").Value = 251
$ws.Range("G43").Value = 2.7
$ws.Range("I43").Value = 2.5
$ws.Range("L43").Value = 3.2
$ws.Range("M43").Value = 1.03
$ws.Range("O43").Value = 1.27
$ws.Range("T43").Value = 1.25
$ws.Range("W43").Value = 1.8
$ws.Range("X43").Value = 1.91
$ws.Range("Y43").Value = 8.5
$ws.Range("M44").Value = 1.03
$ws.Range("O44").Value = 1.22
$ws.Range("T44").Value = 1.33
$ws.Range("AB45").Value = 23
$ws.Range("AL45").Value = 11
$ws.Range("AM45").Value = 29
$ws.Range("AN45").Value = 23
$ws.Range("G45").Value = 2.5
$ws.Range("I45").Value = 2.75
$ws.Range("J45").Value = 3.2
$ws.Range("M45").Value = 1.04
$ws.Range("O45").Value = 1.27
$ws.Range("T45").Value = 1.25
$ws.Range("Y45").Value = 8
